# Fruta / hortaliza, semanal
# Insert a new week's worth of data (3 rows) at the top of the "Palta" price
# table on sheet1, pushing the existing historical rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 758, shifting rows 758:806
# down to 761:809 (and everything in between along with them).
$ws.Rows("758:760").Insert()

# Populate the 3 newly inserted rows with this week's "Especial / Primera /
# Segunda" entries for "Terminal La Palmera de La Serena" (Provincia de
# Limarí), dated 44706.

# Row 758 - Especial
$ws.Range("A758").Value = 8
$ws.Range("B758").Value = "Terminal La Palmera de La Serena"
$ws.Range("C758").Value = "Coquimbo"
$ws.Range("D758").Value = 44706
$ws.Range("E758").Value = 4
$ws.Range("F758").Value = "Fruta"
$ws.Range("G758").Value = 100106
$ws.Range("H758").Value = "Oleaginosos"
$ws.Range("I758").Value = 100106002
$ws.Range("J758").Value = "Palta"
$ws.Range("K758").Value = "Hass"
$ws.Range("L758").Value = "Especial"
$ws.Range("M758").Value = 200
$ws.Range("N758").Value = 3100
$ws.Range("O758").Value = 3200
$ws.Range("P758").Value = 3150
$ws.Range("Q758").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R758").Value = "Provincia de Limarí"
$ws.Range("S758").Value = 3150
$ws.Range("T758").Value = 1

# Row 759 - Primera
$ws.Range("A759").Value = 8
$ws.Range("B759").Value = "Terminal La Palmera de La Serena"
$ws.Range("C759").Value = "Coquimbo"
$ws.Range("D759").Value = 44706
$ws.Range("E759").Value = 4
$ws.Range("F759").Value = "Fruta"
$ws.Range("G759").Value = 100106
$ws.Range("H759").Value = "Oleaginosos"
$ws.Range("I759").Value = 100106002
$ws.Range("J759").Value = "Palta"
$ws.Range("K759").Value = "Hass"
$ws.Range("L759").Value = "Primera"
$ws.Range("M759").Value = 300
$ws.Range("N759").Value = 2900
$ws.Range("O759").Value = 3000
$ws.Range("P759").Value = 2950
$ws.Range("Q759").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R759").Value = "Provincia de Limarí"
$ws.Range("S759").Value = 2950
$ws.Range("T759").Value = 1

# Row 760 - Segunda
$ws.Range("A760").Value = 8
$ws.Range("B760").Value = "Terminal La Palmera de La Serena"
$ws.Range("C760").Value = "Coquimbo"
$ws.Range("D760").Value = 44706
$ws.Range("E760").Value = 4
$ws.Range("F760").Value = "Fruta"
$ws.Range("G760").Value = 100106
$ws.Range("H760").Value = "Oleaginosos"
$ws.Range("I760").Value = 100106002
$ws.Range("J760").Value = "Palta"
$ws.Range("K760").Value = "Hass"
$ws.Range("L760").Value = "Segunda"
$ws.Range("M760").Value = 240
$ws.Range("N760").Value = 2700
$ws.Range("O760").Value = 2800
$ws.Range("P760").Value = 2750
$ws.Range("Q760").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R760").Value = "Provincia de Limarí"
$ws.Range("S760").Value = 2750
$ws.Range("T760").Value = 1

Write-Host "Edit applied"
